# Relabel the food "type" column (B) to the new normalized/merged category
# names. Column A ("code") values are untouched - only the human readable
# labels in column B change, per the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "american,steakhouse"
$ws.Range("B3").Value  = "asian,vietnamese,thai"
$ws.Range("B4").Value  = "bakeries,delis"
$ws.Range("B5").Value  = "barbeque"
$ws.Range("B6").Value  = "breakfast and brunch"
$ws.Range("B7").Value  = "burgers"
$ws.Range("B8").Value  = "coffee and tea"
$ws.Range("B9").Value  = "french"
$ws.Range("B10").Value = "ice cream and frozen yogurt"
$ws.Range("B11").Value = "irish"
$ws.Range("B12").Value = "italian"
$ws.Range("B13").Value = "juice bars and smoothies,cafes"
$ws.Range("B14").Value = "mexican,tex mex,southwestern,tapas,small plates"
$ws.Range("B15").Value = "mediterranean"
$ws.Range("B16").Value = "pizza"
$ws.Range("B17").Value = "pretzels"
$ws.Range("B18").Value = "salad,soup,sandwiches"
$ws.Range("B19").Value = "wine bars,wine and spirits,beer,irish pub,sports bars"

# Update the view state: the user had scrolled down and selected E19
# before saving.
$ws.Range("E19").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
